# sm_car_data_Linkage_Link5_S2LAF.xlsx
# v2p14. Compatible with MF-Swift v2212, updated hardpoints.
#
# This script reproduces the data/formatting edits captured in the
# authoritative OOXML diff: updated Endstop (xMax/xMin) hardpoint values
# on both the front ("...f") and rear ("...r") sheets, the corresponding
# number-format tidy-up on those rows, and the sheet tab colour change.

$wb = $excel.ActiveWorkbook

$wsFront = $wb.Worksheets.Item(1)   # L5StoLAF_Sedan_HambaLG_f
$wsRear  = $wb.Worksheets.Item(2)   # L5StoLAF_Sedan_HambaLG_r

# ---------------------------------------------------------------------
# Tab colour: theme "Accent5" (index 8) tinted -0.249977111117893
# (was theme "Accent4" (index 7) tinted 0.59999389629810485)
# ---------------------------------------------------------------------
$wsFront.Tab.Color = 11957550
$wsRear.Tab.Color  = 11957550

# ---------------------------------------------------------------------
# Front sheet (L5StoLAF_Sedan_HambaLG_f) hardpoint updates
# ---------------------------------------------------------------------

# Row 25 - TrackRod sInboard
$wsFront.Range("F25").Value = 0.15379999999999999
$wsFront.Range("G25").Value = 0.65
$wsFront.Range("H25").Value = 0.24

# Row 26 - TrackRod sOutboard (F26 unchanged)
$wsFront.Range("G26").Value = 0.91
$wsFront.Range("H26").Value = 0.23

# Row 28 - Shock sTop (F28 unchanged)
$wsFront.Range("G28:H29").NumberFormat = "0.00"
$wsFront.Range("G28").Value = 0.62
$wsFront.Range("H28").Value = 0.65

# Row 29 - Shock sBottom (F29 unchanged)
$wsFront.Range("G29").Value = 0.85
$wsFront.Range("H29").Value = 0.19

# ---------------------------------------------------------------------
# Rear sheet (L5StoLAF_Sedan_HambaLG_r) hardpoint updates
# ---------------------------------------------------------------------

# Row 25 - TrackRod sInboard
$wsRear.Range("F25").Value = 0.13
$wsRear.Range("G25").Value = 0.65
$wsRear.Range("H25").Value = 0.24

# Row 26 - TrackRod sOutboard
$wsRear.Range("F26").Value = 0.13
$wsRear.Range("G26").Value = 0.91
$wsRear.Range("H26").Value = 0.23

# Row 28 - Shock sTop
$wsRear.Range("F28:H29").NumberFormat = "0.00"
$wsRear.Range("F28").Value = 0.002655714285714287
$wsRear.Range("G28").Value = 0.62
$wsRear.Range("H28").Value = 0.65

# Row 29 - Shock sBottom
$wsRear.Range("F29").Value = -0.05516642857142858
$wsRear.Range("G29").Value = 0.85
$wsRear.Range("H29").Value = 0.19
